$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2:C5").Value = (Get-Date -Year 2023 -Month 10 -Day 8 -Hour 0 -Minute 0 -Second 0).Date
